$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: insert a new column before DO (the 01-oct. column) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DO1").EntireColumn.Insert()
$ws1.Range("DO1").Value = "10-nov"
$ws1.Range("DO2:DO25").Value = "-"

# --- "Gaz" sheet: append two new daily rows ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Cells.Item(147, 1).Value = "'2025-11-08"
$ws2.Cells.Item(147, 1).Style = "Normal"
$ws2.Cells.Item(147, 2).Value = 29.755
$ws2.Cells.Item(148, 1).Value = "'2025-11-09"
$ws2.Cells.Item(148, 1).Style = "Normal"
$ws2.Cells.Item(148, 2).Value = 29.755

# --- "CO2" sheet: append two new daily rows ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Cells.Item(147, 1).Value = "'2025-11-08"
$ws3.Cells.Item(147, 1).Style = "Normal"
$ws3.Cells.Item(147, 2).Value = 79.36
$ws3.Cells.Item(148, 1).Value = "'2025-11-09"
$ws3.Cells.Item(148, 1).Style = "Normal"
$ws3.Cells.Item(148, 2).Value = 79.36
